# Update the Income Statement "Overview" sheet:
#  - drop the oldest reporting period (1396/12) and its publish-date column,
#    shifting all period/date/figure columns one to the left
#  - append the newest period (1401/12) with its publish date and figures
#  - refresh the restated prior-period figures that came back from the source
#  - the D15 "-" placeholder becomes a real 0 figure along the way

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: financial period headers ---------------------------------------
$ws.Range("D8").Value2 = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value2 = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value2 = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value2 = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value2 = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ----------------------------------------------------
$ws.Range("D9").Value2 = "1399-05-09 (11)"
$ws.Range("E9").Value2 = "1400-05-07 (12)"
$ws.Range("F9").Value2 = "1401-05-09 (9)"
$ws.Range("G9").Value2 = "1402-02-30 (8)"
$ws.Range("H9").Value2 = "1402-02-30 (2)"

# --- Rows 11-27: income-statement figures ------------------------------------
$values = @{
    11 = @(234757238, 391458791, 774036745, 1456266489, 1607038537)
    12 = @(-124803148, -225458652, -406302318, -701949250, -1036494655)
    13 = @(109954090, 166000139, 367734427, 754317239, 570543882)
    14 = @(-7692537, -11938176, -18679133, -34327741, -53307287)
    15 = @(0, 0, 0, 0, 0)
    16 = @(17076969, -1475350, 6356145, 12254290, 23322833)
    17 = @(119338522, 152586613, 355411439, 732243788, 540559428)
    18 = @(-10253754, -14710422, -30243308, -24945630, -37181720)
    19 = @(23084419, 31122266, 90009897, 184678606, 100883934)
    20 = @(132169187, 168998457, 415178028, 891976764, 604261642)
    21 = @(-9864015, -21038649, -32403017, -51182605, -56957617)
    22 = @(122305172, 147959808, 382775011, 840794159, 547304025)
    23 = @(0, 0, 0, 0, 0)
    24 = @(122305172, 147959808, 382775011, 840794159, 547304025)
    25 = @(585, 505, 1306, 2870, 684)
    26 = @(209000000, 293000000, 293000000, 293000000, 800000000)
    27 = @(153, 185, 478, 1051, 684)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("D" + $row).Value2 = $rowValues[0]
    $ws.Range("E" + $row).Value2 = $rowValues[1]
    $ws.Range("F" + $row).Value2 = $rowValues[2]
    $ws.Range("G" + $row).Value2 = $rowValues[3]
    $ws.Range("H" + $row).Value2 = $rowValues[4]
}
